{"js": "const replacements = [\n  [\"41-40=1\", \"3+17=20\"],\n  [\"86+0=86\", \"58-11=47\"],\n  [\"42-28=14\", \"74-61=13\"],\n  [\"59-51=8\", \"93-51=42\"],\n  [\"23+4=27\", \"2+12=14\"],\n  [\"22+17=39\", \"32+46=78\"],\n  [\"76-12=64\", \"50-20=30\"],\n  [\"11+19=30\", \"93-39=54\"],\n  [\"59-3=56\", \"98-52=46\"],\n  [\"32+29=61\", \"71-47=24\"],\n  [\"47+31=78\", \"53+12=65\"],\n  [\"6+24=30\", \"79-32=47\"],\n  [\"50+28=78\", \"76-28=48\"],\n  [\"35-17=18\", \"3+29=32\"],\n  [\"66+33=99\", \"84-47=37\"],\n  [\"32+31=63\", \"45+3=48\"],\n  [\"47-19=28\", \"6+34=40\"],\n  [\"44-24=20\", \"25+34=59\"],\n  [\"54-51=3\", \"54-10=44\"],\n  [\"9+46=55\", \"78+13=91\"],\n  [\"70-36=34\", \"43+36=79\"],\n  [\"25-7=18\", \"78-18=60\"],\n  [\"53-35=18\", \"38+11=49\"],\n  [\"53+34=87\", \"14+34=48\"],\n  [\"6+40=46\", \"86-43=43\"],\n  [\"3+52=55\", \"61+22=83\"],\n  [\"82-42=40\", \"80-24=56\"],\n  [\"87-29=58\", \"77-27=50\"],\n  [\"23+45=68\", \"0+30=30\"],\n  [\"56-24=32\", \"72+11=83\"],\n  [\"27+18=45\", \"87+11=98\"],\n  [\"40+19=59\", \"49-11=38\"],\n  [\"86-67=19\", \"53-49=4\"],\n  [\"5+43=48\", \"0+66=66\"],\n  [\"47+29=76\", \"4+16=20\"],\n  [\"7+51=58\", \"78-39=39\"],\n  [\"97-40=57\", \"5+8=13\"],\n  [\"44-31=13\", \"38+11=49\"],\n  [\"96-57=39\", \"90+8=98\"],\n  [\"95-41=54\", \"20-4=16\"],\n  [\"30+18=48\", \"33-4=29\"],\n  [\"7+8=15\", \"50-3=47\"],\n  [\"97-94=3\", \"11+35=46\"],\n  [\"76-25=51\", \"28-27=1\"],\n  [\"20+19=39\", \"27-4=23\"],\n  [\"53-0=53\", \"23+73=96\"],\n  [\"63-2=61\", \"9+16=25\"],\n  [\"67-26=41\", \"73-14=59\"],\n  [\"9+42=51\", \"51-44=7\"],\n  [\"34+26=60\", \"42+16=58\"],\n  [\"38-27=11\", \"59-32=27\"],\n  [\"79-36=43\", \"15+34=49\"],\n  [\"44-18=26\", \"12+33=45\"],\n  [\"11+44=55\", \"46+22=68\"],\n  [\"85-8=77\", \"5+15=20\"],\n  [\"21+30=51\", \"24-21=3\"],\n  [\"95-27=68\", \"63-31=32\"],\n  [\"9+37=46\", \"50-5=45\"],\n  [\"20+23=43\", \"14+8=22\"],\n  [\"86-48=38\", \"63+24=87\"],\n  [\"28+1=29\", \"33+17=50\"],\n  [\"92-33=59\", \"47-41=6\"],\n  [\"20+8=28\", \"73-66=7\"],\n  [\"71+10=81\", \"6+38=44\"],\n  [\"45+23=68\", \"67-56=11\"],\n  [\"70-58=12\", \"51+3=54\"],\n  [\"23-0=23\", \"40+45=85\"],\n  [\"37+25=62\", \"20-8=12\"],\n  [\"40-10=30\", \"61-11=50\"],\n  [\"76+15=91\", \"25+29=54\"],\n  [\"25+57=82\", \"67+3=70\"],\n  [\"96-39=57\", \"80-7=73\"],\n  [\"26+37=63\", \"78-54=24\"],\n  [\"53-24=29\", \"27+24=51\"],\n  [\"1+49=50\", \"49-3=46\"],\n  [\"68+10=78\", \"85-36=49\"],\n  [\"58+8=66\", \"3+84=87\"],\n  [\"19-3=16\", \"63-62=1\"],\n  [\"42+29=71\", \"9+65=74\"],\n  [\"30+66=96\", \"93+6=99\"],\n  [\"92-45=47\", \"93-47=46\"],\n  [\"74-7=67\", \"95-85=10\"],\n  [\"80+10=90\", \"47-4=43\"],\n  [\"22+44=66\", \"28+67=95\"],\n  [\"4+55=59\", \"77-47=30\"],\n  [\"53+23=76\", \"97-13=84\"],\n  [\"81-22=59\", \"18+62=80\"],\n  [\"85-51=34\", \"44-41=3\"],\n  [\"31+55=86\", \"66-64=2\"],\n  [\"61+31=92\", \"5+86=91\"],\n  [\"43+34=77\", \"18+26=44\"],\n  [\"89-82=7\", \"70-32=38\"],\n  [\"85+12=97\", \"65-23=42\"],\n  [\"9+8=17\", \"31-5=26\"],\n  [\"99-81=18\", \"62+15=77\"],\n  [\"94-70=24\", \"79-64=15\"],\n  [\"50-35=15\", \"59-58=1\"],\n  [\"96-43=53\", \"38-36=2\"],\n  [\"42+56=98\", \"40+39=79\"],\n  [\"81-23=58\", \"22+15=37\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}", "ps1": "# Replace each old equation result with its updated value, cell by cell,\n# using Word's Find/Replace (each old string is unique in the document).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('41-40=1', '3+17=20'),\n    @('86+0=86', '58-11=47'),\n    @('42-28=14', '74-61=13'),\n    @('59-51=8', '93-51=42'),\n    @('23+4=27', '2+12=14'),\n    @('22+17=39', '32+46=78'),\n    @('76-12=64', '50-20=30'),\n    @('11+19=30', '93-39=54'),\n    @('59-3=56', '98-52=46'),\n    @('32+29=61', '71-47=24'),\n    @('47+31=78', '53+12=65'),\n    @('6+24=30', '79-32=47'),\n    @('50+28=78', '76-28=48'),\n    @('35-17=18', '3+29=32'),\n    @('66+33=99', '84-47=37'),\n    @('32+31=63', '45+3=48'),\n    @('47-19=28', '6+34=40'),\n    @('44-24=20', '25+34=59'),\n    @('54-51=3', '54-10=44'),\n    @('9+46=55', '78+13=91'),\n    @('70-36=34', '43+36=79'),\n    @('25-7=18', '78-18=60'),\n    @('53-35=18', '38+11=49'),\n    @('53+34=87', '14+34=48'),\n    @('6+40=46', '86-43=43'),\n    @('3+52=55', '61+22=83'),\n    @('82-42=40', '80-24=56'),\n    @('87-29=58', '77-27=50'),\n    @('23+45=68', '0+30=30'),\n    @('56-24=32', '72+11=83'),\n    @('27+18=45', '87+11=98'),\n    @('40+19=59', '49-11=38'),\n    @('86-67=19', '53-49=4'),\n    @('5+43=48', '0+66=66'),\n    @('47+29=76', '4+16=20'),\n    @('7+51=58', '78-39=39'),\n    @('97-40=57', '5+8=13'),\n    @('44-31=13', '38+11=49'),\n    @('96-57=39', '90+8=98'),\n    @('95-41=54', '20-4=16'),\n    @('30+18=48', '33-4=29'),\n    @('7+8=15', '50-3=47'),\n    @('97-94=3', '11+35=46'),\n    @('76-25=51', '28-27=1'),\n    @('20+19=39', '27-4=23'),\n    @('53-0=53', '23+73=96'),\n    @('63-2=61', '9+16=25'),\n    @('67-26=41', '73-14=59'),\n    @('9+42=51', '51-44=7'),\n    @('34+26=60', '42+16=58'),\n    @('38-27=11', '59-32=27'),\n    @('79-36=43', '15+34=49'),\n    @('44-18=26', '12+33=45'),\n    @('11+44=55', '46+22=68'),\n    @('85-8=77', '5+15=20'),\n    @('21+30=51', '24-21=3'),\n    @('95-27=68', '63-31=32'),\n    @('9+37=46', '50-5=45'),\n    @('20+23=43', '14+8=22'),\n    @('86-48=38', '63+24=87'),\n    @('28+1=29', '33+17=50'),\n    @('92-33=59', '47-41=6'),\n    @('20+8=28', '73-66=7'),\n    @('71+10=81', '6+38=44'),\n    @('45+23=68', '67-56=11'),\n    @('70-58=12', '51+3=54'),\n    @('23-0=23', '40+45=85'),\n    @('37+25=62', '20-8=12'),\n    @('40-10=30', '61-11=50'),\n    @('76+15=91', '25+29=54'),\n    @('25+57=82', '67+3=70'),\n    @('96-39=57', '80-7=73'),\n    @('26+37=63', '78-54=24'),\n    @('53-24=29', '27+24=51'),\n    @('1+49=50', '49-3=46'),\n    @('68+10=78', '85-36=49'),\n    @('58+8=66', '3+84=87'),\n    @('19-3=16', '63-62=1'),\n    @('42+29=71', '9+65=74'),\n    @('30+66=96', '93+6=99'),\n    @('92-45=47', '93-47=46'),\n    @('74-7=67', '95-85=10'),\n    @('80+10=90', '47-4=43'),\n    @('22+44=66', '28+67=95'),\n    @('4+55=59', '77-47=30'),\n    @('53+23=76', '97-13=84'),\n    @('81-22=59', '18+62=80'),\n    @('85-51=34', '44-41=3'),\n    @('31+55=86', '66-64=2'),\n    @('61+31=92', '5+86=91'),\n    @('43+34=77', '18+26=44'),\n    @('89-82=7', '70-32=38'),\n    @('85+12=97', '65-23=42'),\n    @('9+8=17', '31-5=26'),\n    @('99-81=18', '62+15=77'),\n    @('94-70=24', '79-64=15'),\n    @('50-35=15', '59-58=1'),\n    @('96-43=53', '38-36=2'),\n    @('42+56=98', '40+39=79'),\n    @('81-23=58', '22+15=37'),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}\n"}
